$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28 (shifts existing rows 28-66 down to 29-67,
# carrying their formatting/data along for the ride).
$ws.Rows(28).Insert()

# Populate the newly inserted row 28 with the new weekly price record
# (Orégano, Vega Central Mapocho de Santiago) for date 44720.
$ws.Range("A28").Value = 9
$ws.Range("B28").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C28").Value = "Metropolitana"
$ws.Range("D28").Value = 44720
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = 100112029
$ws.Range("G28").Value = "Orégano"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 15000
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = 15000
$ws.Range("N28").Value = "$/docena de atados"
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 5000
$ws.Range("Q28").Value = 3
$ws.Range("R28").Value = "Hortaliza"
